# Insert a new data row at row 186 (pushing existing rows 186-234 down to 187-235)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("186:186").Insert()

$ws.Cells.Item(186, 1).Value = 7
$ws.Cells.Item(186, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(186, 3).Value = "Ñuble"
$ws.Cells.Item(186, 4).Value = 44754
$ws.Cells.Item(186, 5).Value = 16
$ws.Cells.Item(186, 6).Value = 100112043
$ws.Cells.Item(186, 7).Value = "Pepino ensalada"
$ws.Cells.Item(186, 8).Value = "Sin especificar"
$ws.Cells.Item(186, 9).Value = "Primera"
$ws.Cells.Item(186, 10).Value = 100
$ws.Cells.Item(186, 11).Value = 19000
$ws.Cells.Item(186, 12).Value = 20000
$ws.Cells.Item(186, 13).Value = 19500
$ws.Cells.Item(186, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(186, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(186, 16).Value = 325
$ws.Cells.Item(186, 17).Value = 60
$ws.Cells.Item(186, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same number format as the rest of column D
$ws.Cells.Item(186, 4).NumberFormat = $ws.Cells.Item(187, 4).NumberFormat
